# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> hashtable of row number -> new F-column value.
$updates = @{
    "展览"     = @{ 2 = 677; 3 = 28; 5 = 1889; 6 = 42; 7 = 3283; 8 = 460; 9 = 786 }
    "全部类型" = @{ 2 = 677; 3 = 28; 6 = 1889; 7 = 42; 8 = 3283; 9 = 460; 10 = 786 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
